# "Future Plan page is created"
#
# 1) Fix the title of slide 9 ("Weakness of Our Project"): the original
#    title is split across two runs ("Weakness of " + "Our Project") with
#    a trailing endParaRPr. Re-typing the full text collapses it back into
#    a single run and drops the stray endParaRPr.
$p = $ppt.ActivePresentation

$weaknessSlide = $p.Slides.Item(9)
$titleRange = $weaknessSlide.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = ""
$titleRange = $weaknessSlide.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Weakness of Our Project"

# 2) Append a new "Future Plan" slide at the end of the deck, using the
#    same "Title and Content" layout as the rest of the presentation.
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)

$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Future Plan"

$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = " "
$body.ParagraphFormat.Bullet.Type = 0

$rulerLevel = $newSlide.Shapes.Item(2).TextFrame.Ruler.Levels.Item(1)
$rulerLevel.LeftMargin = 0
$rulerLevel.FirstMargin = 0
